$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (04-30-2015, Thursday): fix the late/time-out/overtime figures.
$ws.Range("C10").Value = "13:04:00"
$ws.Range("D10").Value = "15:51:00"
$ws.Range("I10").Value = 1.0

# Row 14 (05-04-2015, Monday): add a time-in, drop overtime to half day,
# and re-colour the row the same as the other "has attendance" rows (row 7/10/18).
$ws.Range("C14").Value = "13:40:00"
$ws.Range("I14").Value = 0.5
$ws.Range("A14:P14").Interior.Color = $ws.Range("A7").Interior.Color

# Row 16 (05-06-2015, Wednesday): add time-in/out, undertime and adjust overtime.
$ws.Range("C16").Value = "13:58:00"
$ws.Range("D16").Value = "17:13:00"
$ws.Range("F16").Value = 1.5
$ws.Range("I16").Value = 0.5
$ws.Range("A16:P16").Interior.Color = $ws.Range("A7").Interior.Color

# Row 17 (05-07-2015, Thursday): add time-in/out, undertime and adjust overtime.
$ws.Range("C17").Value = "15:05:00"
$ws.Range("D17").Value = "17:55:00"
$ws.Range("F17").Value = 0.75
$ws.Range("I17").Value = 0.5
$ws.Range("A17:P17").Interior.Color = $ws.Range("A7").Interior.Color
